# Adds the new Artisan Command "pidWeights(<beta>,<gamma>)" to the
# "Commands" worksheet, right after the existing "p-i-d(<p>,<i>,<d>)"
# command row, pushing every row below it down by one (adjustSV, pidSV,
# ... through the trailing WebSocket/Stepper/RC command rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# "p-i-d(<p>,<i>,<d>)" / "sets the p-i-d parameters of the PID" live on
# row 108; the new row is inserted directly below it (new row 109),
# which shifts the old row 109 ("adjustSV(<float>)") and everything
# after it down by one row, matching the target layout.
$ws.Rows.Item(109).Insert()

$ws.Range("B109").Value2 = "pidWeights(<beta>,<gamma>)"
$ws.Range("C109").Value2 = "sets the beta and gamma parameters of the PID"

# Update the visible selection/scroll position to match the edited area.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 97
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C109").Select()
